$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Re-solved "adjustable cells" (Solver changed D5:E5) ---
$ws.Range("D5").Value = 150
$ws.Range("E5").Value = 20.000000000000007

# --- Updated inputs ---
$ws.Range("E13").Value = 1000      # Publico Estimado TV budget split
$ws.Range("D17").Value = 6000      # nº anúncios no Rádio (Total row)
$ws.Range("E17").Value = 14000     # nº anúncios na TV (Total row)

# --- Re-apply formatting on the "Total" row so it matches the other data
#     rows instead of the old (fill-applied) style ---
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D16:F16").Copy()
$ws.Range("D17:F17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update the Solver "solver_opt" defined name so it targets G9 ---
$wb.Names.Item("solver_opt").RefersTo = "=Sheet1!`$G`$9"

# --- Move the active selection, as last left by the user ---
$ws.Range("E10").Select()

$wb.Save()
